$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# hw6 actual (column D) - mirror the claimed values in column C
$ws.Range("D6").Value = 12
$ws.Range("D10").Value = 4
$ws.Range("D11").Value = 4
$ws.Range("D12").Value = 2
$ws.Range("D13").Value = 4
$ws.Range("D14").Value = 4
$ws.Range("D15").Value = 5

# hw7 claimed (column E) - mirror the possible values in column B
$ws.Range("E26").Value = 4
$ws.Range("E28").Value = 5
$ws.Range("E29").Value = 2
$ws.Range("E30").Value = 2

# Update sheet view: zoom, top-left cell, active cell selection
$ws.Application.ActiveWindow.Zoom = 145
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I23").Select()
